$d = $word.ActiveDocument

# 1) Update the letter date.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "September 19, 2025`r") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# 2) Split the mailing address "2960 Feliz Rd, Santa Clara CA 95051" into
#    two separate paragraphs: "2960 Feliz Rd" and "Santa Clara, CA 95051".
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "2960 Feliz Rd, Santa Clara CA 95051`r") {
        $p.Range.Text = "2960 Feliz Rd"
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = "Santa Clara, CA 95051"
        break
    }
}

# 3) Remove the empty "No Spacing" paragraph directly after the
#    "... Board of Directors" line.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Townhomes at Nuevo Homeowners Association Board of Directors`r") {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}
